$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "item" sheet: insert a new row at 582 for item id 1901
#    "Spurning Ribbon" (shifts every row below it down by one).
# ---------------------------------------------------------------------
$itemWs = $wb.Worksheets.Item("item")

$itemWs.Rows.Item(582).Insert()
$itemWs.Cells.Item(582, 1).Value = 1901
$itemWs.Cells.Item(582, 2).Value = "Spurning Ribbon"

# Match the formatting used by the row above (B119 carries the same
# "no explicit left-alignment" style that the new row needs) instead of
# the left-aligned style Excel would otherwise copy down from row 583.
$itemWs.Cells.Item(119, 2).Copy()
$itemWs.Cells.Item(582, 2).PasteSpecial(-4122)
$itemWs.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) "item_inventory" sheet: row 420 previously had id 0 pointing at
#    "Spurning Ribbon" - now it should point at item id 1901.
# ---------------------------------------------------------------------
$invWs = $wb.Worksheets.Item("item_inventory")
$invWs.Cells.Item(420, 1).Value = 1901
$invWs.Cells.Item(420, 2).Value = "Spurning Ribbon"

# ---------------------------------------------------------------------
# 3) View/selection state: the author ended the session with
#    item_inventory active (selection on E423), while "item" and
#    "treasure_states" retain plain (non-active) selections.
# ---------------------------------------------------------------------
$itemWs.Range("B582").Select()

$treasureWs = $wb.Worksheets.Item("treasure_states")
$treasureWs.Range("B167").Select()

$invWs.Activate()
$invWs.Range("E423").Select()
